$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the TDL_EC_49_12S row (old row 2); remaining rows shift up by one
$ws.Rows.Item(2).Delete()

# Update Richness (B), Shannon (C) and Simpson (D) values for all sample rows
$ws.Cells.Item(2, 2).Value = 12
$ws.Cells.Item(2, 3).Value = 0.65
$ws.Cells.Item(2, 4).Value = 0.314
$ws.Cells.Item(3, 2).Value = 12
$ws.Cells.Item(3, 3).Value = 0.843
$ws.Cells.Item(3, 4).Value = 0.482
$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(4, 3).Value = 0.478
$ws.Cells.Item(4, 4).Value = 0.235
$ws.Cells.Item(5, 2).Value = 20
$ws.Cells.Item(5, 3).Value = 0.988
$ws.Cells.Item(5, 4).Value = 0.458
$ws.Cells.Item(6, 2).Value = 20
$ws.Cells.Item(6, 3).Value = 0.894
$ws.Cells.Item(6, 4).Value = 0.475
$ws.Cells.Item(7, 2).Value = 13
$ws.Cells.Item(7, 3).Value = 0.551
$ws.Cells.Item(7, 4).Value = 0.276
$ws.Cells.Item(8, 2).Value = 13
$ws.Cells.Item(8, 3).Value = 0.147
$ws.Cells.Item(8, 4).Value = 0.044
$ws.Cells.Item(9, 2).Value = 16
$ws.Cells.Item(9, 3).Value = 0.603
$ws.Cells.Item(9, 4).Value = 0.211
$ws.Cells.Item(10, 2).Value = 22
$ws.Cells.Item(10, 3).Value = 1.934
$ws.Cells.Item(10, 4).Value = 0.769
$ws.Cells.Item(11, 2).Value = 18
$ws.Cells.Item(11, 3).Value = 1.465
$ws.Cells.Item(11, 4).Value = 0.641
$ws.Cells.Item(12, 2).Value = 20
$ws.Cells.Item(12, 3).Value = 2.162
$ws.Cells.Item(12, 4).Value = 0.8110000000000001
$ws.Cells.Item(13, 2).Value = 23
$ws.Cells.Item(13, 3).Value = 2.315
$ws.Cells.Item(13, 4).Value = 0.857
$ws.Cells.Item(14, 2).Value = 16
$ws.Cells.Item(14, 3).Value = 1.837
$ws.Cells.Item(14, 4).Value = 0.746
$ws.Cells.Item(15, 2).Value = 14
$ws.Cells.Item(15, 3).Value = 2.124
$ws.Cells.Item(15, 4).Value = 0.839
$ws.Cells.Item(16, 2).Value = 22
$ws.Cells.Item(16, 3).Value = 1.863
$ws.Cells.Item(16, 4).Value = 0.786
$ws.Cells.Item(17, 2).Value = 22
$ws.Cells.Item(17, 3).Value = 1.666
$ws.Cells.Item(17, 4).Value = 0.735
$ws.Cells.Item(18, 2).Value = 26
$ws.Cells.Item(18, 3).Value = 1.386
$ws.Cells.Item(18, 4).Value = 0.633
$ws.Cells.Item(19, 2).Value = 22
$ws.Cells.Item(19, 3).Value = 1.544
$ws.Cells.Item(19, 4).Value = 0.709
$ws.Cells.Item(20, 2).Value = 32
$ws.Cells.Item(20, 3).Value = 1.86
$ws.Cells.Item(20, 4).Value = 0.784
$ws.Cells.Item(21, 2).Value = 18
$ws.Cells.Item(21, 3).Value = 1.724
$ws.Cells.Item(21, 4).Value = 0.746
$ws.Cells.Item(22, 2).Value = 17
$ws.Cells.Item(22, 3).Value = 1.553
$ws.Cells.Item(22, 4).Value = 0.709
$ws.Cells.Item(23, 2).Value = 22
$ws.Cells.Item(23, 3).Value = 1.669
$ws.Cells.Item(23, 4).Value = 0.757
$ws.Cells.Item(24, 2).Value = 27
$ws.Cells.Item(24, 3).Value = 2.104
$ws.Cells.Item(24, 4).Value = 0.8149999999999999
$ws.Cells.Item(25, 2).Value = 22
$ws.Cells.Item(25, 3).Value = 0.773
$ws.Cells.Item(25, 4).Value = 0.294
$ws.Cells.Item(26, 2).Value = 26
$ws.Cells.Item(26, 3).Value = 1.888
$ws.Cells.Item(26, 4).Value = 0.722
$ws.Cells.Item(27, 2).Value = 24
$ws.Cells.Item(27, 3).Value = 1.714
$ws.Cells.Item(27, 4).Value = 0.728
$ws.Cells.Item(28, 2).Value = 29
$ws.Cells.Item(28, 3).Value = 1.25
$ws.Cells.Item(28, 4).Value = 0.544
$ws.Cells.Item(29, 2).Value = 22
$ws.Cells.Item(29, 3).Value = 1.613
$ws.Cells.Item(29, 4).Value = 0.74
$ws.Cells.Item(30, 2).Value = 26
$ws.Cells.Item(30, 3).Value = 1.938
$ws.Cells.Item(30, 4).Value = 0.738
$ws.Cells.Item(31, 2).Value = 29
$ws.Cells.Item(31, 3).Value = 2.034
$ws.Cells.Item(31, 4).Value = 0.751
$ws.Cells.Item(32, 2).Value = 43
$ws.Cells.Item(32, 3).Value = 1.989
$ws.Cells.Item(32, 4).Value = 0.744
$ws.Cells.Item(33, 2).Value = 16
$ws.Cells.Item(33, 3).Value = 1.182
$ws.Cells.Item(33, 4).Value = 0.544
$ws.Cells.Item(34, 2).Value = 23
$ws.Cells.Item(34, 3).Value = 1.938
$ws.Cells.Item(34, 4).Value = 0.785
$ws.Cells.Item(35, 2).Value = 21
$ws.Cells.Item(35, 3).Value = 1.68
$ws.Cells.Item(35, 4).Value = 0.719
$ws.Cells.Item(36, 2).Value = 26
$ws.Cells.Item(36, 3).Value = 2.077
$ws.Cells.Item(36, 4).Value = 0.8
$ws.Cells.Item(37, 2).Value = 20
$ws.Cells.Item(37, 3).Value = 1.853
$ws.Cells.Item(37, 4).Value = 0.759
$ws.Cells.Item(38, 2).Value = 21
$ws.Cells.Item(38, 3).Value = 1.706
$ws.Cells.Item(38, 4).Value = 0.723
$ws.Cells.Item(39, 2).Value = 25
$ws.Cells.Item(39, 3).Value = 1.715
$ws.Cells.Item(39, 4).Value = 0.732
$ws.Cells.Item(40, 2).Value = 19
$ws.Cells.Item(40, 3).Value = 1.299
$ws.Cells.Item(40, 4).Value = 0.62
$ws.Cells.Item(41, 2).Value = 43
$ws.Cells.Item(41, 3).Value = 2.306
$ws.Cells.Item(41, 4).Value = 0.779
$ws.Cells.Item(42, 2).Value = 39
$ws.Cells.Item(42, 3).Value = 2.612
$ws.Cells.Item(42, 4).Value = 0.885
$ws.Cells.Item(43, 2).Value = 13
$ws.Cells.Item(43, 3).Value = 1.118
$ws.Cells.Item(43, 4).Value = 0.511
$ws.Cells.Item(44, 2).Value = 16
$ws.Cells.Item(44, 3).Value = 1.014
$ws.Cells.Item(44, 4).Value = 0.435
$ws.Cells.Item(45, 2).Value = 12
$ws.Cells.Item(45, 3).Value = 0.432
$ws.Cells.Item(45, 4).Value = 0.178
$ws.Cells.Item(46, 2).Value = 15
$ws.Cells.Item(46, 3).Value = 0.534
$ws.Cells.Item(46, 4).Value = 0.24
$ws.Cells.Item(47, 2).Value = 19
$ws.Cells.Item(47, 3).Value = 0.851
$ws.Cells.Item(47, 4).Value = 0.527
$ws.Cells.Item(48, 2).Value = 16
$ws.Cells.Item(48, 3).Value = 1.256
$ws.Cells.Item(48, 4).Value = 0.671
